$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new topic columns X1:AC1 (columns 24-29)
$ws.Cells.Item(1, 24).Value2 = "Israel Gaza War"
$ws.Cells.Item(1, 25).Value2 = "GPT-4"
$ws.Cells.Item(1, 26).Value2 = "MS Dhoni"
$ws.Cells.Item(1, 27).Value2 = "Moye Moye meme"
$ws.Cells.Item(1, 28).Value2 = "Al Nassr"
$ws.Cells.Item(1, 29).Value2 = "Inter Miami"

# Data rows 2-51 for the 6 new topic columns (X:AC / columns 24-29)
$data = @{
  2 = @(0, 1, 11, 0, 37, 0)
  3 = @(0, 2, 12, 1, 16, 0)
  4 = @(0, 0, 17, 0, 100, 0)
  5 = @(0, 0, 16, 1, 24, 0)
  6 = @(0, 0, 14, 0, 14, 0)
  7 = @(0, 1, 13, 0, 13, 0)
  8 = @(0, 0, 12, 0, 14, 0)
  9 = @(0, 0, 14, 1, 15, 0)
  10 = @(1, 1, 13, 0, 14, 0)
  11 = @(0, 1, 11, 1, 15, 0)
  12 = @(0, 100, 12, 0, 14, 0)
  13 = @(0, 30, 16, 0, 7, 0)
  14 = @(0, 14, 34, 1, 4, 0)
  15 = @(1, 16, 38, 0, 12, 0)
  16 = @(0, 15, 42, 3, 10, 0)
  17 = @(0, 11, 36, 1, 9, 0)
  18 = @(0, 10, 39, 1, 12, 0)
  19 = @(0, 10, 40, 0, 4, 0)
  20 = @(1, 11, 43, 0, 6, 0)
  21 = @(0, 10, 38, 0, 6, 0)
  22 = @(1, 9, 50, 0, 8, 0)
  23 = @(0, 8, 100, 0, 8, 0)
  24 = @(0, 9, 45, 0, 4, 28)
  25 = @(0, 8, 29, 0, 3, 10)
  26 = @(0, 8, 23, 0, 2, 4)
  27 = @(1, 8, 22, 0, 3, 5)
  28 = @(1, 7, 72, 0, 2, 6)
  29 = @(0, 10, 21, 0, 3, 9)
  30 = @(0, 11, 20, 0, 9, 43)
  31 = @(0, 9, 14, 0, 21, 46)
  32 = @(0, 8, 13, 0, 19, 62)
  33 = @(0, 6, 12, 1, 34, 92)
  34 = @(0, 7, 19, 0, 20, 68)
  35 = @(0, 10, 15, 0, 20, 100)
  36 = @(0, 7, 19, 0, 18, 55)
  37 = @(0, 8, 18, 0, 9, 42)
  38 = @(1, 9, 22, 0, 9, 27)
  39 = @(0, 6, 17, 0, 19, 38)
  40 = @(0, 7, 15, 0, 17, 33)
  41 = @(13, 6, 18, 1, 15, 22)
  42 = @(100, 7, 18, 7, 5, 12)
  43 = @(80, 8, 17, 1, 8, 9)
  44 = @(51, 7, 20, 3, 16, 10)
  45 = @(51, 7, 22, 4, 15, 4)
  46 = @(36, 9, 23, 5, 13, 12)
  47 = @(27, 7, 31, 12, 6, 5)
  48 = @(16, 8, 51, 88, 11, 4)
  49 = @(10, 6, 25, 100, 16, 4)
  50 = @(13, 8, 17, 39, 12, 4)
  51 = @(10, 8, 17, 24, 14, 4)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($c = 0; $c -lt $vals.Count; $c++) {
    $ws.Cells.Item($row, 24 + $c).Value2 = $vals[$c]
  }
}

